$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value is numeric-looking but must keep a trailing
# zero (e.g. "1.030"), which Excel would otherwise normalize to "1.03" when
# assigned through .Value. Force them to Text format first so the literal
# string is preserved exactly.
$textCells = @('D6', 'D18', 'D19', 'D22', 'D25', 'D28', 'D43')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.693.89'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '1.848.84'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('D4').Value = '1.031'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '321.02'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('D6').Value = '1.030'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').Value = '0.4387'
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').Value = '0.3787'
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').Value = '0.8819'
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').Value = '21.53'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.880.14'
$ws.Range('E12').Value = '  +1.82%  '
$ws.Range('D13').Value = '5.496'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = '6.688'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '0.07161'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '84.88'
$ws.Range('E16').Value = '  +2.60%  '
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '0.000009070'
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').Value = '1.030'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').Value = '15.46'
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('D21').Value = '27.720.60'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').Value = '5.280'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('D24').Value = '2.103.98'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('D25').Value = '2.050'
$ws.Range('E25').Value = '  +6.37%  '
$ws.Range('D26').Value = '158.86'
$ws.Range('E26').Value = '  +0.70%  '
$ws.Range('D27').Value = '18.67'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').Value = '1.990'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('D30').Value = '117.53'
$ws.Range('E30').Value = '  +1.34%  '
$ws.Range('D31').Value = '0.09066'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').Value = '0.7721'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('D34').Value = '3.004'
$ws.Range('E34').Value = '  +4.45%  '
$ws.Range('D35').Value = '4.552'
$ws.Range('E35').Value = '  +1.19%  '
$ws.Range('D36').Value = '1.031'
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('D37').Value = '1.149'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').Value = '0.01974'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').Value = '0.05257'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').Value = '2.845'
$ws.Range('E40').Value = '  +2.16%  '
$ws.Range('D41').Value = '0.5175'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').Value = '0.1669'
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').Value = '6.860'
$ws.Range('E43').Value = '  +2.81%  '
$ws.Range('E44').Value = '  +1.74%  '
$ws.Range('D45').Value = '110.16'
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').Value = '10.67'
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('D48').Value = '0.06566'
$ws.Range('E48').Value = '  +3.25%  '
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D50').Value = '0.4691'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('E51').Value = '  -0.68%  '
